$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to G16 (matches the sheetView selection change in the diff)
$ws.Range("G16").Select()

# Row 9: was a header-like row (style matching row 4), becomes a normal data row
# (style matching rows 5-8 / 10-14). Copy the formatting from an existing data
# row so the cell-format (fill/number-format) matches exactly, then update
# the values: status changes from "Em análise" to "Desenvolvido" and a
# delivery date of 10/02/2013 (serial 41315) is added in column D.
$ws.Range("A5:C5").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C9").Value = "Desenvolvido"
$ws.Range("D9").Value = "02/10/2013"

# Rows 10-14: delivery date changes from 07/12/2013 (serial 41615) to 10/02/2013 (serial 41315)
$ws.Range("D10").Value = "02/10/2013"
$ws.Range("D11").Value = "02/10/2013"
$ws.Range("D12").Value = "02/10/2013"
$ws.Range("D13").Value = "02/10/2013"
$ws.Range("D14").Value = "02/10/2013"
